# Apply updated crypto price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.717.15"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.888.46"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'247.85"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.4739"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.2928"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'0.06533"
$ws.Range("D10").Value = "'22.00"
$ws.Range("D11").Value = "'0.07801"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "1.891.10"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "'0.7361"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'5.246"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "'284.65"
$ws.Range("E16").Value = "  +4.27%  "
$ws.Range("D17").Value = "30.697.86"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "'13.22"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "'0.000007535"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D21").Value = "2.139.41"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "'5.340"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").Value = "'0.9995"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'6.258"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'9.233"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'164.56"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'18.93"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'1.924"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").Value = "'0.09733"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "'4.304"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'4.187"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "'2.721"
$ws.Range("D38").Value = "'0.01906"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("D39").Value = "'2.802"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "'6.402"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").Value = "'76.04"
$ws.Range("E41").Value = "  +7.19%  "
$ws.Range("D42").Value = "'2.010"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").Value = "'0.4261"
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").Value = "'0.9993"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'0.8354"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'101.68"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "'9.541"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "'35.68"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").Value = "'7.030"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'919.13"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'0.05751"
$ws.Range("E51").Value = "  +2.14%  "
